$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-6 (old fiscal years 16-17 through 20-21) so the remaining
# data (previously rows 7-25) shifts up to become rows 2-20, matching the
# new algorithmically-generated distribution output.
$ws.Range("A2:D6").EntireRow.Delete()
